$wb = $excel.ActiveWorkbook

# --- Sheet "CIRULO PANEL": assign component values ---
$ws1 = $wb.Worksheets.Item("CIRULO PANEL")

# C13 capacitor (row 8) gets its Value set to 100n
$ws1.Range("E8").Value = "100n"

# R40 resistor (row 16) gets its Value set to 220 (numeric)
$ws1.Range("E16").Value = 220

# restore selection on the sheet
$null = $ws1.Range("B3").Select()

# --- Sheet "_HISTORY": log the change ---
$ws2 = $wb.Worksheets.Item("_HISTORY")

$ws2.Range("A5").Value = 2
$ws2.Range("B5").Value = 44950
$ws2.Range("B5").NumberFormat = $ws2.Range("B4").NumberFormat
$ws2.Range("C5").Value = "JRC"
$ws2.Range("D5").Value = "Valor de Condensadores"

$null = $ws2.Range("C6").Select()

$null = $wb.Save()
